$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(1112, 1, 889, 1009, 8, 3, 977, 24315, 13, 107, 6062)
    3  = @(256, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    4  = @(173, 0, 71, 76, 5, 0, 658, 1690, 1, 0, 0)
    5  = @(1512, 2, 76, 83, 5, 0, 1679, 27754, 1, 0, 0)
    6  = @(1543, 0, 256, 271, 15, 2, 2598, 75143, 8, 0, 0)
    7  = @(118, 0, 21, 26, 4, 0, 1538, 4286, 0, 0, 0)
    8  = @(1019, 0, 505, 547, 4, 0, 110, 15267, 7, 39, 1077)
    9  = @(255, 0, 7, 7, 0, 0, 0, 5714, 0, 0, 0)
    10 = @(1316, 0, 498, 519, 18, 3, 2258, 13472, 5, 0, 0)
    11 = @(1765, 0, 130, 132, 2, 0, 1524, 14205, 0, 0, 0)
    12 = @(1374, 0, 98, 131, 33, 0, 5377, 6717, 2, 0, 0)
    13 = @(850, 0, 219, 222, 1, 2, 1250, 13201, 8, 0, 0)
    14 = @(257, 0, 13, 13, 0, 0, 0, 6154, 1, 0, 0)
    15 = @(416, 0, 52, 51, 0, 0, 0, 30784, 0, 0, 0)
    16 = @(177, 0, 4, 4, 0, 0, 0, 0, 0, 0, 0)
    17 = @(1751, 0, 19, 20, 1, 0, 1429, 22143, 0, 0, 0)
    18 = @(4175, 0, 294, 619, 15, 3, 1015, 26791, 7, 306, 18527)
}

$cols = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12)  # B..L

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item($row, $cols[$i]).Value = $values[$i]
    }
}
